$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '58.620.16'
$ws.Range('E2').Value = '  -4.10%  '
$ws.Range('D3').Value = '2.558.89'
$ws.Range('D4').Value = "'1.00"
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').Value = "'508.94"
$ws.Range('E5').Value = '  -4.65%  '
$ws.Range('D6').Value = "'146.00"
$ws.Range('E6').Value = '  -6.74%  '
$ws.Range('D7').Value = "'0.998"
$ws.Range('E7').Value = '  +0.05%  '
$ws.Range('D8').Value = "'0.569"
$ws.Range('E8').Value = '  -3.78%  '
$ws.Range('D9').Value = '2.570.73'
$ws.Range('E9').Value = '  -3.88%  '
$ws.Range('D10').Value = "'6.21"
$ws.Range('E10').Value = '  -5.65%  '
$ws.Range('E11').Value = '  -6.29%  '
$ws.Range('D12').Value = "'0.335"
$ws.Range('E12').Value = '  -5.05%  '
$ws.Range('D14').Value = '3.009.87'
$ws.Range('E14').Value = '  -3.42%  '
$ws.Range('D15').Value = '58.574.97'
$ws.Range('E15').Value = '  -4.10%  '
$ws.Range('D16').Value = "'20.97"
$ws.Range('E16').Value = '  -5.18%  '
$ws.Range('D17').Value = "'0.0000137"
$ws.Range('E17').Value = '  -5.29%  '
$ws.Range('D18').Value = '2.562.65'
$ws.Range('E18').Value = '  -3.90%  '
$ws.Range('B19').Value = 'BitcoinCash'
$ws.Range('C19').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D19').Value = "'346.68"
$ws.Range('E19').Value = '  -2.53%  '
$ws.Range('B20').Value = 'Polkadot'
$ws.Range('C20').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D20').Value = "'4.54"
$ws.Range('E20').Value = '  -5.04%  '
$ws.Range('D21').Value = "'10.22"
$ws.Range('E21').Value = '  -4.47%  '
$ws.Range('B22').Value = 'Uniswap'
$ws.Range('C22').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D22').Value = "'6.00"
$ws.Range('E22').Value = '  -4.28%  '
$ws.Range('B23').Value = 'Dai'
$ws.Range('C23').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D23').Value = "'1.00"
$ws.Range('E23').Value = '  +0.21%  '
$ws.Range('D24').Value = "'60.63"
$ws.Range('E24').Value = '  -1.72%  '
$ws.Range('D25').Value = "'0.413"
$ws.Range('E25').Value = '  -4.44%  '
$ws.Range('D26').Value = "'0.996"
$ws.Range('E26').Value = '  -0.19%  '
$ws.Range('D27').Value = "'0.160"
$ws.Range('E27').Value = '  -5.12%  '
$ws.Range('D28').Value = '2.658.98'
$ws.Range('E28').Value = '  -3.97%  '
$ws.Range('D29').Value = '0.0₃0800'
$ws.Range('E29').Value = '  -7.18%  '
$ws.Range('D30').Value = "'7.01"
$ws.Range('E30').Value = '  -5.24%  '
$ws.Range('D31').Value = "'0.998"
$ws.Range('E31').Value = '  -0.12%  '
$ws.Range('D32').Value = "'5.97"
$ws.Range('E32').Value = '  -3.84%  '
$ws.Range('B33').Value = 'EthereumClassic'
$ws.Range('C33').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D33').Value = "'18.64"
$ws.Range('E33').Value = '  -4.90%  '
$ws.Range('B34').Value = 'Monero'
$ws.Range('C34').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D34').Value = "'149.47"
$ws.Range('E34').Value = '  -0.43%  '
$ws.Range('E35').Value = '  -6.29%  '
$ws.Range('D36').Value = "'3.96"
$ws.Range('E36').Value = '  -4.52%  '
$ws.Range('D37').Value = "'0.897"
$ws.Range('E37').Value = '  +0.84%  '
$ws.Range('E38').Value = '  -6.72%  '
$ws.Range('D39').Value = "'0.842"
$ws.Range('E39').Value = '  -8.76%  '
$ws.Range('D40').Value = "'36.01"
$ws.Range('E40').Value = '  -2.30%  '
$ws.Range('D41').Value = "'1.40"
$ws.Range('E41').Value = '  -6.58%  '
$ws.Range('D42').Value = "'287.11"
$ws.Range('E42').Value = '  -6.99%  '
$ws.Range('D43').Value = "'3.56"
$ws.Range('E43').Value = '  -7.04%  '
$ws.Range('D44').Value = "'0.0997"
$ws.Range('E44').Value = '  -2.62%  '
$ws.Range('D45').Value = "'0.996"
$ws.Range('E45').Value = '  -0.05%  '
$ws.Range('D46').Value = "'0.605"
$ws.Range('E46').Value = '  -6.88%  '
$ws.Range('D47').Value = "'0.0536"
$ws.Range('E47').Value = '  -5.16%  '
$ws.Range('D48').Value = "'19.02"
$ws.Range('E48').Value = '  -5.33%  '
$ws.Range('B49').Value = 'VeChain'
$ws.Range('C49').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D49').Value = "'0.0229"
$ws.Range('E49').Value = '  -4.54%  '
$ws.Range('B50').Value = 'WhiteBITCoin'
$ws.Range('C50').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D50').Value = "'10.25"
$ws.Range('E50').Value = '  -1.05%  '
$ws.Range('D51').Value = "'4.68"
$ws.Range('E51').Value = '  -7.61%  '

foreach ($addr in @('D4','D5','D6','D7','D8','D10','D12','D16','D17','D19','D20','D21','D22','D23','D24','D25','D26','D27','D30','D31','D32','D33','D34','D36','D37','D39','D40','D41','D42','D43','D44','D45','D46','D47','D48','D49','D50','D51')) {
    $ws.Range($addr).Style = 'Normal'
}
